$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.530.26"
$ws.Range("E2").Value = "'  +1.37%  "
$ws.Range("E3").Value = "'  +1.16%  "
$ws.Range("D4").Value = "'1.0000"
$ws.Range("E4").Value = "'  -0.16%  "
$ws.Range("D5").Value = "'233.65"
$ws.Range("E5").Value = "'  +2.41%  "
$ws.Range("E6").Value = "'  -0.21%  "
$ws.Range("D7").Value = "'0.4745"
$ws.Range("E7").Value = "'  +3.19%  "
$ws.Range("E8").Value = "'  +2.50%  "
$ws.Range("D9").Value = "'0.06321"
$ws.Range("E9").Value = "'  +1.76%  "
$ws.Range("D10").Value = "'17.73"
$ws.Range("E10").Value = "'  +11.50%  "
$ws.Range("D11").Value = "'1.822.59"
$ws.Range("E11").Value = "'  -0.74%  "
$ws.Range("D12").Value = "'0.07454"
$ws.Range("E12").Value = "'  +1.66%  "
$ws.Range("D13").Value = "'4.976"
$ws.Range("E13").Value = "'  +2.13%  "
$ws.Range("D14").Value = "'84.71"
$ws.Range("E14").Value = "'  +2.27%  "
$ws.Range("D15").Value = "'0.6276"
$ws.Range("E15").Value = "'  +2.06%  "
$ws.Range("D16").Value = "'30.492.89"
$ws.Range("D17").Value = "'246.21"
$ws.Range("E17").Value = "'  +9.26%  "
$ws.Range("E18").Value = "'  -0.29%  "
$ws.Range("E19").Value = "'  +3.60%  "
$ws.Range("D20").Value = "'0.000007340"
$ws.Range("E20").Value = "'  +2.11%  "
$ws.Range("E21").Value = "'  +0.03%  "
$ws.Range("D22").Value = "'4.947"
$ws.Range("E22").Value = "'  +2.33%  "
$ws.Range("E23").Value = "'  +1.77%  "
$ws.Range("D24").Value = "'9.143"
$ws.Range("E24").Value = "'  +0.87%  "
$ws.Range("D25").Value = "'162.72"
$ws.Range("E25").Value = "'  -1.45%  "
$ws.Range("E26").Value = "'  +2.50%  "
$ws.Range("D27").Value = "'1.878"
$ws.Range("E27").Value = "'  +1.96%  "
$ws.Range("D28").Value = "'0.1026"
$ws.Range("E28").Value = "'  +1.73%  "
$ws.Range("D29").Value = "'1.352"
$ws.Range("E29").Value = "'  -1.72%  "
$ws.Range("D30").Value = "'4.027"
$ws.Range("E30").Value = "'  -0.58%  "
$ws.Range("D31").Value = "'3.842"
$ws.Range("E31").Value = "'  +2.76%  "
$ws.Range("D32").Value = "'0.04849"
$ws.Range("E32").Value = "'  +1.24%  "
$ws.Range("E33").Value = "'  +0.41%  "
$ws.Range("D34").Value = "'0.7029"
$ws.Range("E34").Value = "'  +1.90%  "
$ws.Range("D35").Value = "'2.697"
$ws.Range("E35").Value = "'  -0.21%  "
$ws.Range("D36").Value = "'0.01900"
$ws.Range("E36").Value = "'  +5.13%  "
$ws.Range("D37").Value = "'2.685"
$ws.Range("E37").Value = "'  +3.37%  "
$ws.Range("D38").Value = "'2.002"
$ws.Range("E38").Value = "'  +4.77%  "
$ws.Range("D39").Value = "'0.8748"
$ws.Range("E39").Value = "'  -0.64%  "
$ws.Range("D40").Value = "'106.67"
$ws.Range("E40").Value = "'  +3.75%  "
$ws.Range("D41").Value = "'1.000"
$ws.Range("E41").Value = "'  +0.34%  "
$ws.Range("E42").Value = "'  +1.98%  "
$ws.Range("D43").Value = "'0.4066"
$ws.Range("E43").Value = "'  +2.37%  "
$ws.Range("D44").Value = "'7.231"
$ws.Range("E44").Value = "'  +5.54%  "
$ws.Range("D45").Value = "'62.85"
$ws.Range("E45").Value = "'  +6.82%  "
$ws.Range("D46").Value = "'0.1204"
$ws.Range("E46").Value = "'  +2.41%  "
$ws.Range("D47").Value = "'33.68"
$ws.Range("E47").Value = "'  +4.00%  "
$ws.Range("D48").Value = "'8.523"
$ws.Range("E48").Value = "'  +1.43%  "
$ws.Range("E49").Value = "'  -0.09%  "
$ws.Range("D50").Value = "'1.354"
$ws.Range("E50").Value = "'  -0.10%  "
$ws.Range("D51").Value = "'0.3695"
$ws.Range("E51").Value = "'  +2.59%  "
